$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 (ECs) ---
$ws.Range("G2").Value = 0.4859026666666666
$ws.Range("H2").Value = 1.457708
$ws.Range("M2").Value = 4.662797333333334
$ws.Range("N2").Value = 13.988392
$ws.Range("O2").Value = 0.7324994586787992
$ws.Range("P2").Value = 0.7324994586787993
$ws.Range("Q2").Value = 2.265665658392889
$ws.Range("R2").Value = 20.390990925536
$ws.Range("S2").Value = 0.7324994586787992
$ws.Range("T2").Value = 0.7324994586787993

# --- Update existing row 3 (FAPs) ---
$ws.Range("G3").Value = 0.4859026666666666
$ws.Range("H3").Value = 1.457708
$ws.Range("O3").Value = 0.1045598489170565
$ws.Range("P3").Value = 0.1045598489170565
$ws.Range("Q3").Value = 0.3234100122959999
$ws.Range("R3").Value = 2.910690110664
$ws.Range("S3").Value = 0.1045598489170565
$ws.Range("T3").Value = 0.1045598489170565

# --- Update existing row 4 (MuSCs) ---
$ws.Range("G4").Value = 0.4859026666666666
$ws.Range("H4").Value = 1.457708
$ws.Range("M4").Value = 0.7894166666666665
$ws.Range("N4").Value = 2.36825
$ws.Range("O4").Value = 0.1240129561007488
$ws.Range("P4").Value = 0.1240129561007488
$ws.Range("Q4").Value = 0.3835796634444444
$ws.Range("R4").Value = 3.452216970999999
$ws.Range("S4").Value = 0.1240129561007488
$ws.Range("T4").Value = 0.1240129561007488

# --- New row 5 (Neutrophils) ---
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Fgf5"
$ws.Range("C5").Value = "Fgfr3"
$ws.Range("D5").Value = "Neutrophils"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.4859026666666666
$ws.Range("H5").Value = 1.457708
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.08057833333333334
$ws.Range("N5").Value = 0.241735
$ws.Range("O5").Value = 0.01265840681643176
$ws.Range("P5").Value = 0.01265840681643176
$ws.Range("Q5").Value = 0.03915322704222222
$ws.Range("R5").Value = 0.35237904338
$ws.Range("S5").Value = 0.01265840681643176
$ws.Range("T5").Value = 0.01265840681643176

# --- New row 6 (Resolving-Mac) ---
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Fgf5"
$ws.Range("C6").Value = "Fgfr3"
$ws.Range("D6").Value = "Resolving-Mac"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.4859026666666666
$ws.Range("H6").Value = 1.457708
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 1
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.16722
$ws.Range("N6").Value = 0.50166
$ws.Range("O6").Value = 0.02626932948696365
$ws.Range("P6").Value = 0.02626932948696365
$ws.Range("Q6").Value = 0.08125264392000001
$ws.Range("R6").Value = 0.73127379528
$ws.Range("S6").Value = 0.02626932948696365
$ws.Range("T6").Value = 0.02626932948696365
